$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E (reviews_count) - shifts everything after it one column to the left
$ws.Range("E:E").Delete()
